$d = $word.ActiveDocument
$d.Content.Find.Execute(
    "My GXUST student ID is:201700408066, My SCU student ID is 23347492. ",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 2
)
Write-Output $d.Content.Text
